# Fruta / hortaliza, semanal
# The weekly refresh cyclically rotates the D (Fecha) and L:T (Calidad .. Kg/unidad)
# fields among rows 2, 3, 6, 7, 8 of the sheet:
#   row2 <- old row7, row3 <- old row8, row6 <- old row3, row7 <- old row6, row8 <- old row2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return @{
        D = $ws.Cells.Item($row, 4).Value2
        L = $ws.Cells.Item($row, 12).Value2
        M = $ws.Cells.Item($row, 13).Value2
        N = $ws.Cells.Item($row, 14).Value2
        O = $ws.Cells.Item($row, 15).Value2
        P = $ws.Cells.Item($row, 16).Value2
        Q = $ws.Cells.Item($row, 17).Value2
        R = $ws.Cells.Item($row, 18).Value2
        S = $ws.Cells.Item($row, 19).Value2
        T = $ws.Cells.Item($row, 20).Value2
    }
}

function Set-RowData($row, $data) {
    $ws.Cells.Item($row, 4).Value = $data.D
    $ws.Cells.Item($row, 12).Value = $data.L
    $ws.Cells.Item($row, 13).Value = $data.M
    $ws.Cells.Item($row, 14).Value = $data.N
    $ws.Cells.Item($row, 15).Value = $data.O
    $ws.Cells.Item($row, 16).Value = $data.P
    $ws.Cells.Item($row, 17).Value = $data.Q
    $ws.Cells.Item($row, 18).Value = $data.R
    $ws.Cells.Item($row, 19).Value = $data.S
    $ws.Cells.Item($row, 20).Value = $data.T
}

$row2 = Get-RowData 2
$row3 = Get-RowData 3
$row6 = Get-RowData 6
$row7 = Get-RowData 7
$row8 = Get-RowData 8

Set-RowData 2 $row7
Set-RowData 3 $row8
Set-RowData 6 $row3
Set-RowData 7 $row6
Set-RowData 8 $row2

Write-Output "rotation applied"
